$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: TbCommon / Common / TRUE / 名词解释表.xlsx
# Reuse the "good" (green) cell style from an existing styled row by copying
# formats first, then writing the values so no new style entries are created.
$ws.Range("B5").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("C5").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("E5").Copy()
$ws.Range("E8").PasteSpecial(-4122)

$ws.Range("B8").Value = "TbCommon"
$ws.Range("C8").Value = "Common"
$ws.Range("D8").Value = $true
$ws.Range("E8").Value = "名词解释表.xlsx"

# Row 9: TbBuffInfo / BuffInfo / TRUE / BuffInfo.xlsx / BuffName
$ws.Range("B5").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("C5").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("E5").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("E5").Copy()
$ws.Range("F9").PasteSpecial(-4122)

$ws.Range("B9").Value = "TbBuffInfo"
$ws.Range("C9").Value = "BuffInfo"
$ws.Range("D9").Value = $true
$ws.Range("E9").Value = "BuffInfo.xlsx"
$ws.Range("F9").Value = "BuffName"

$ws.Application.CutCopyMode = $false

# Match the saved selection state recorded in the target workbook.
$ws.Range("F9").Select() | Out-Null
